$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3577
$ws.Range("F5").Value = 3577
$ws.Range("G5").Value = 75
$ws.Range("F7").Value = 5107
$ws.Range("F8").Value = 5107
$ws.Range("F9").Value = 525
$ws.Range("F10").Value = 358
$ws.Range("F12").Value = 693
$ws.Range("F14").Value = 88
$ws.Range("F15").Value = 33
$ws.Range("F16").Value = 700
$ws.Range("F17").Value = 318
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 91
$ws.Range("F23").Value = 4913
$ws.Range("F24").Value = 4913
$ws.Range("F28").Value = 6042
$ws.Range("F32").Value = 344
$ws.Range("F33").Value = 712
$ws.Range("F34").Value = 4445
$ws.Range("F35").Value = 319
$ws.Range("F36").Value = 123
$ws.Range("F38").Value = 1020
$ws.Range("F40").Value = 24
$ws.Range("F42").Value = 872
$ws.Range("F43").Value = 994
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 12
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 227
$ws.Range("F3").Value = 1119
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 227
$ws.Range("F4").Value = 1119
$ws.Range("F7").Value = 3577
$ws.Range("F8").Value = 3577
$ws.Range("G8").Value = 75
$ws.Range("F10").Value = 5107
$ws.Range("F11").Value = 5107
$ws.Range("F12").Value = 525
$ws.Range("F13").Value = 358
$ws.Range("F15").Value = 693
$ws.Range("F17").Value = 88
$ws.Range("F18").Value = 33
$ws.Range("F19").Value = 700
$ws.Range("F20").Value = 318
$ws.Range("F21").Value = 35
$ws.Range("F23").Value = 91
$ws.Range("F27").Value = 4913
$ws.Range("F28").Value = 4913
$ws.Range("F32").Value = 6042
$ws.Range("F36").Value = 344
$ws.Range("F37").Value = 712
$ws.Range("F38").Value = 4445
$ws.Range("F39").Value = 319
$ws.Range("F41").Value = 123
$ws.Range("F43").Value = 1020
$ws.Range("F45").Value = 24
$ws.Range("F47").Value = 872
$ws.Range("F48").Value = 994
$ws.Range("F49").Value = 12
